$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'41.157.72"
$ws.Range("E2").Value = "'  -3.50%  "
$ws.Range("D3").Value = "'2.446.39"
$ws.Range("E3").Value = "'  -3.40%  "
$ws.Range("E4").Value = "'  -0.01%  "
$ws.Range("D5").Value = "'309.73"
$ws.Range("E5").Value = "'  +0.51%  "
$ws.Range("D6").Value = "'93.10"
$ws.Range("E6").Value = "'  -7.36%  "
$ws.Range("D7").Value = "'0.550"
$ws.Range("E7").Value = "'  -3.58%  "
$ws.Range("E8").Value = "'  +0.08%  "
$ws.Range("D9").Value = "'0.502"
$ws.Range("E9").Value = "'  -4.76%  "
$ws.Range("D10").Value = "'33.23"
$ws.Range("E10").Value = "'  -7.24%  "
$ws.Range("D11").Value = "'0.0779"
$ws.Range("E11").Value = "'  -3.17%  "
$ws.Range("D12").Value = "'0.108"
$ws.Range("E12").Value = "'  -0.47%  "
$ws.Range("D13").Value = "'6.94"
$ws.Range("E13").Value = "'  -5.25%  "
$ws.Range("D14").Value = "'2.816.44"
$ws.Range("E14").Value = "'  -3.90%  "
$ws.Range("D15").Value = "'2.447.40"
$ws.Range("E15").Value = "'  -5.10%  "
$ws.Range("D16").Value = "'14.38"
$ws.Range("E16").Value = "'  -9.34%  "
$ws.Range("D17").Value = "'0.783"
$ws.Range("E17").Value = "'  -3.38%  "
$ws.Range("D18").Value = "'41.121.85"
$ws.Range("E18").Value = "'  -3.55%  "
$ws.Range("D19").Value = "'6.32"
$ws.Range("E19").Value = "'  -6.42%  "
$ws.Range("D20").Value = "'0.0₃0909"
$ws.Range("E20").Value = "'  -4.53%  "
$ws.Range("D21").Value = "'11.47"
$ws.Range("E21").Value = "'  -6.26%  "
$ws.Range("D22").Value = "'67.45"
$ws.Range("E22").Value = "'  -2.72%  "
$ws.Range("D23").Value = "'236.36"
$ws.Range("E23").Value = "'  -3.02%  "
$ws.Range("E24").Value = "'  -4.29%  "
$ws.Range("D25").Value = "'1.93"
$ws.Range("E25").Value = "'  -5.74%  "
$ws.Range("E26").Value = "'  +0.09%  "
$ws.Range("D27").Value = "'24.45"
$ws.Range("E27").Value = "'  -5.99%  "
$ws.Range("E28").Value = "'  -5.57%  "
$ws.Range("D29").Value = "'9.65"
$ws.Range("E29").Value = "'  -4.93%  "
$ws.Range("D30").Value = "'35.92"
$ws.Range("E30").Value = "'  -8.29%  "
$ws.Range("D31").Value = "'151.91"
$ws.Range("E31").Value = "'  -2.30%  "
$ws.Range("D32").Value = "'5.58"
$ws.Range("E32").Value = "'  -3.59%  "
$ws.Range("D33").Value = "'2.61"
$ws.Range("E33").Value = "'  -0.76%  "
$ws.Range("D36").Value = "'3.01"
$ws.Range("E36").Value = "'  -4.71%  "
$ws.Range("D37").Value = "'17.12"
$ws.Range("E37").Value = "'  -7.26%  "
$ws.Range("D38").Value = "'1.88"
$ws.Range("E38").Value = "'  -7.39%  "
$ws.Range("D39").Value = "'0.104"
$ws.Range("E39").Value = "'  -7.66%  "
$ws.Range("D40").Value = "'0.114"
$ws.Range("E40").Value = "'  -4.28%  "
$ws.Range("D41").Value = "'4.15"
$ws.Range("E41").Value = "'  -3.15%  "
$ws.Range("D42").Value = "'21.16"
$ws.Range("E42").Value = "'  -4.84%  "
$ws.Range("D43").Value = "'1.00"
$ws.Range("E43").Value = "'  +0.09%  "
$ws.Range("D44").Value = "'1.967.23"
$ws.Range("E44").Value = "'  +0.04%  "
$ws.Range("D45").Value = "'0.0283"
$ws.Range("E45").Value = "'  -4.96%  "
$ws.Range("D46").Value = "'3.03"
$ws.Range("E46").Value = "'  -7.59%  "
$ws.Range("D47").Value = "'8.71"
$ws.Range("E47").Value = "'  -1.95%  "
$ws.Range("D48").Value = "'76.70"
$ws.Range("E48").Value = "'  -5.10%  "
$ws.Range("D49").Value = "'96.82"
$ws.Range("E49").Value = "'  -3.98%  "
$ws.Range("D50").Value = "'68.72"
$ws.Range("E50").Value = "'  -5.13%  "
$ws.Range("D51").Value = "'0.179"
$ws.Range("E51").Value = "'  -6.53%  "

# Rows 34-35: Hedera and ApeXProtocol swap places (with updated price/volume)
$ws.Range("B34").Value = "'Hedera"
$ws.Range("C34").Value = "'https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D34").Value = "'0.0750"
$ws.Range("E34").Value = "'  -5.25%  "
$ws.Range("B35").Value = "'ApeXProtocol"
$ws.Range("C35").Value = "'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D35").Value = "'2.52"
$ws.Range("E35").Value = "'  -7.86%  "
